$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

$ws.Range("B2").Value = "2 "
$ws.Range("C2").Value = " Contabil "
$ws.Range("D2").Value = " SIA"
$ws.Range("E2").Value = "3 "
$ws.Range("F2").Value = " Contabil "
$ws.Range("G2").Value = " SAVA"

$ws.Rows("3:5").Delete()

$ws.Range("B13").Select()
